$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:6").Copy()
$ws.Rows("7:7").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)
$ws.Range("B7:E7").Clear()
$ws.Range("A7").Value = "Distribution"
$ws.Range("A7").Select()
